# Remove "删除" column from detail view
# This removes the two header cells for is_locked_lbl (G1) and
# is_enabled_lbl (H1), shifting the remaining columns (order_by, rem)
# left so they become G1/H1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1:H1").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftToLeft)
